$wb = $excel.ActiveWorkbook

# --- Overview sheet: status text updated (affects column widths implicitly) ---
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("E2").Value = "Handed back: in sync with en-US"
$ovw.Range("F2").Value = "Handed back: in sync with en-US"
$ovw.Range("E3").Value = "Handed back: in sync with en-US"
$ovw.Range("F3").Value = "Handed back: in sync with en-US"
$ovw.Columns.Item(5).ColumnWidth = 29.9777047293527
$ovw.Columns.Item(6).ColumnWidth = 29.9777047293527

# --- zh-cn sheet ---
$zh = $wb.Worksheets.Item("zh-cn")

# Row 2 (ab19d52a...)
$zh.Range("C2").Value = "Handed back: in sync with en-US"
$zh.Hyperlinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7ac1cd02e10acbe463920a359696746559598949/e2e/ab19d52a-605f-4a0d-af24-bedb291dff6c.md", "", "", "ab19d52a-605f-4a0d-af24-bedb291dff6c.md") | Out-Null
$zh.Range("J2").Value = "ab19d52a-605f-4a0d-af24-bedb291dff6c.a0cb47fa219ef42020cdcaa8751b05135353ebba.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-19 04:20:58"

# Row 3 (ea28ad60...)
$zh.Range("C3").Value = "Handed back: in sync with en-US"
$zh.Hyperlinks.Add($zh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7ac1cd02e10acbe463920a359696746559598949/e2e/ea28ad60-4616-4c8c-944a-c3bd9497fa14.md", "", "", "ea28ad60-4616-4c8c-944a-c3bd9497fa14.md") | Out-Null
$zh.Range("J3").Value = "ea28ad60-4616-4c8c-944a-c3bd9497fa14.3a66a22106250585506959241d3046b2244efd80.zh-cn.xlf"
$zh.Range("K3").Value = "2016-08-19 04:20:58"

$zh.Columns.Item(3).ColumnWidth = 29.9777047293527
$zh.Columns.Item(9).ColumnWidth = 40
$zh.Columns.Item(10).ColumnWidth = 40

# --- de-de sheet ---
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = "Handed back: in sync with en-US"
$de.Range("C3").Value = "Handed back: in sync with en-US"

# Row 2 (ab19d52a...)
$de.Hyperlinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7ac1cd02e10acbe463920a359696746559598949/e2e/ab19d52a-605f-4a0d-af24-bedb291dff6c.md", "", "", "ab19d52a-605f-4a0d-af24-bedb291dff6c.md") | Out-Null
$de.Range("J2").Value = "ab19d52a-605f-4a0d-af24-bedb291dff6c.a0cb47fa219ef42020cdcaa8751b05135353ebba.de-de.xlf"
$de.Range("K2").Value = "2016-08-19 04:21:12"

# Row 3 (ea28ad60...)
$de.Hyperlinks.Add($de.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7ac1cd02e10acbe463920a359696746559598949/e2e/ea28ad60-4616-4c8c-944a-c3bd9497fa14.md", "", "", "ea28ad60-4616-4c8c-944a-c3bd9497fa14.md") | Out-Null
$de.Range("J3").Value = "ea28ad60-4616-4c8c-944a-c3bd9497fa14.3a66a22106250585506959241d3046b2244efd80.de-de.xlf"
$de.Range("K3").Value = "2016-08-19 04:21:12"

$de.Columns.Item(3).ColumnWidth = 29.9777047293527
$de.Columns.Item(9).ColumnWidth = 40
$de.Columns.Item(10).ColumnWidth = 40

Write-Host "Edit complete"
